# Applies the "contactdata" diff: row 39 loses its (NaN) email-address
# placeholder and its ID becomes a true number, and four brand-new
# contact rows (40-43) are appended below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39: Alastair Stewart -------------------------------------------
# B39 was stored as text "6"; it should become a real number.
$ws.Cells.Item(39, 2).Value = 6
# F39 (Email Address) had the literal text "NaN"; it should end up blank.
$ws.Cells.Item(39, 6).Value = ""

# --- Helper data for the four new rows ----------------------------------
$newRows = @(
    @{ Row=40; B=1; First="Fraser";  Last="Thorne";  Url="https://www.linkedin.com/in/fraser-thorne-57974510"; Email="";    Company="Edison Group";                  Position="Founder, CEO";                               Connected="2024-11-18 00:00:00"; Gender="Male"; Ethnicity="Caucasian"; Age="50-59"; Industry="Finance/Business" },
    @{ Row=41; B=2; First="Reiss";   Last="Garwood"; Url="https://www.linkedin.com/in/reiss-garwood-635b09b6";  Email="";    Company="eXp Realty";                    Position="Independent property consultant";            Connected="2024-11-18 00:00:00"; Gender="Male"; Ethnicity="Caucasian"; Age="30-39"; Industry="Finance/Business" },
    @{ Row=42; B=3; First="Jonny";   Last="Page";    Url="https://www.linkedin.com/in/jonny-page";              Email="";    Company="Esmée Fairbairn Foundation";    Position="Head of Social and Impact Investment";       Connected="2024-11-15 00:00:00"; Gender="Male"; Ethnicity="Caucasian"; Age="30-39"; Industry="Charity" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    # Column A keeps the same bordered "0" style used throughout the sheet;
    # copy it (value + format) straight from the existing A39 cell.
    $ws.Cells.Item(39, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.First
    $ws.Cells.Item($row, 4).Value = $r.Last
    $ws.Cells.Item($row, 5).Value = $r.Url
    $ws.Cells.Item($row, 6).Value = $r.Email
    $ws.Cells.Item($row, 7).Value = $r.Company
    $ws.Cells.Item($row, 8).Value = $r.Position
    $ws.Cells.Item($row, 9).Value = $r.Connected
    $ws.Cells.Item($row, 10).Value = $r.Gender
    $ws.Cells.Item($row, 11).Value = $r.Ethnicity
    $ws.Cells.Item($row, 12).Value = $r.Age
    $ws.Cells.Item($row, 13).Value = $r.Industry
}

# --- Row 43: Douglas Sloan (B stays text, F keeps its "NaN" placeholder) -
$ws.Cells.Item(39, 1).Copy($ws.Cells.Item(43, 1))
$ws.Cells.Item(43, 2).NumberFormat = "@"
$ws.Cells.Item(43, 2).Value = "4"
$ws.Cells.Item(43, 2).Style = "Normal"
$ws.Cells.Item(43, 3).Value = "Douglas"
$ws.Cells.Item(43, 4).Value = "Sloan"
$ws.Cells.Item(43, 5).Value = "https://www.linkedin.com/in/douglassloan"
$ws.Cells.Item(43, 6).Value = "NaN"
$ws.Cells.Item(43, 7).Value = "ImpactVC"
$ws.Cells.Item(43, 8).Value = "Co-Founder"
$ws.Cells.Item(43, 9).Value = "2024-11-15 00:00:00"
$ws.Cells.Item(43, 10).Value = "Male"
$ws.Cells.Item(43, 11).Value = "Caucasian"
$ws.Cells.Item(43, 12).Value = "30-39"
$ws.Cells.Item(43, 13).Value = "Finance/Business"
